$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '41.392.08'
$ws.Range('E2').Value = '  -3.38%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.474.99'
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '312.60'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '94.15'
$ws.Range('E6').Value = '  -6.57%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.548'
$ws.Range('E7').Value = '  -3.27%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.498'
$ws.Range('E9').Value = '  -4.84%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '33.46'
$ws.Range('E10').Value = '  -5.80%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0781'
$ws.Range('E11').Value = '  -2.94%  '
$ws.Range('E12').Value = '  -0.96%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.98'
$ws.Range('E13').Value = '  -4.57%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.857.19'
$ws.Range('E14').Value = '  -2.68%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.24'
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.468.97'
$ws.Range('E16').Value = '  -4.75%  '
$ws.Range('E17').Value = '  -3.52%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '41.354.89'
$ws.Range('E18').Value = '  -3.46%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.31'
$ws.Range('E19').Value = '  -6.42%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0925'
$ws.Range('E20').Value = '  -2.87%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.24'
$ws.Range('E21').Value = '  -8.81%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '68.71'
$ws.Range('E22').Value = '  -1.74%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '237.08'
$ws.Range('E23').Value = '  -2.63%  '
$ws.Range('E24').Value = '  -4.91%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  -6.15%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '24.02'
$ws.Range('E27').Value = '  -6.69%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.24'
$ws.Range('E28').Value = '  -4.58%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.69'
$ws.Range('E29').Value = '  -4.64%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '36.30'
$ws.Range('E30').Value = '  -5.55%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '151.64'
$ws.Range('E31').Value = '  -4.04%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.46'
$ws.Range('E32').Value = '  -7.55%  '
$ws.Range('E33').Value = '  -3.97%  '
$ws.Range('E34').Value = '  -6.79%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0745'
$ws.Range('E35').Value = '  -6.26%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.07'
$ws.Range('E36').Value = '  -2.96%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '17.35'
$ws.Range('E37').Value = '  -3.44%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.87'
$ws.Range('E38').Value = '  -5.44%  '
$ws.Range('E39').Value = '  -2.96%  '
$ws.Range('E40').Value = '  -8.90%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.21'
$ws.Range('E41').Value = '  +1.43%  '
$ws.Range('E42').Value = '  +0.36%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.991.08'
$ws.Range('E43').Value = '  -0.37%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '19.31'
$ws.Range('E44').Value = '  -11.91%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0285'
$ws.Range('E45').Value = '  -4.81%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.99'
$ws.Range('E46').Value = '  -8.82%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.71'
$ws.Range('E47').Value = '  -4.68%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.720.98'
$ws.Range('E48').Value = '  -2.41%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '69.32'
$ws.Range('E49').Value = '  -4.37%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '96.96'
$ws.Range('E50').Value = '  -4.59%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '74.30'
$ws.Range('E51').Value = '  -6.99%  '
